$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "YYYYQ4" text labels in column A (rows 2-39) with actual
# end-of-year (Q4 / Dec-31) dates, formatted with a new date/time format.
$dateFormat = "YYYY-MM-DD HH:MM:SS"
$startYear = 1987

for ($i = 0; $i -lt 38; $i++) {
    $row = $i + 2
    $year = $startYear + $i
    $d = Get-Date -Year $year -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = $dateFormat
    $cell.Value = $d
}
